# ATT_Bill.xlsx - update June and July 2017 bill entries (sheet "2017")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2017")

# --- Row 8 (2017-06, date 42901): fill in previously-blank billing values,
#     mirroring the pattern already used on row 7. Copy formats first so the
#     newly-populated cells pick up the same currency/alignment styling as
#     row 7 (G8/H8 move from the "39" style to the "38" style).
$ws.Range("G7").Copy() | Out-Null
$ws.Range("G8").PasteSpecial(-4122) | Out-Null
$ws.Range("H7").Copy() | Out-Null
$ws.Range("H8").PasteSpecial(-4122) | Out-Null

$ws.Range("B8").Value = 176.33
$ws.Range("E8").Value = 35.32
$ws.Range("F8").Value = 35.32
$ws.Range("G8").Value = 35.32
$ws.Range("H8").Value = 35.32
$ws.Range("I8").Value = 35.05

# --- Row 9 (2017-07, date 42931): same treatment as row 8.
$ws.Range("G7").Copy() | Out-Null
$ws.Range("G9").PasteSpecial(-4122) | Out-Null
$ws.Range("H7").Copy() | Out-Null
$ws.Range("H9").PasteSpecial(-4122) | Out-Null

$ws.Range("B9").Value = 176.33
$ws.Range("E9").Value = 35.32
$ws.Range("F9").Value = 35.32
$ws.Range("G9").Value = 35.32
$ws.Range("H9").Value = 35.32
$ws.Range("I9").Value = 35.05

# --- Payment-history table (rows 19-24): the $5.32/$5.05 "2x" line items were
#     previously written as 5.28+5.32 / 5.01+5.05; re-expressed as *2 formulas.
$ws.Range("G19").Formula = "=5.32*2"
$ws.Range("G20").Formula = "=5.32*2"
$ws.Range("G21").Formula = "=5.05*2"

# Row 21 gains a new payment-history entry (Balaji, 2017-06-20, $141.20).
$ws.Range("L21").Value = "Balaji"
$ws.Range("M21").Value = 42858
$ws.Range("O21").Value = 141.2

# Rows 22/23: G was a static value, now a formula; F picks up the same
# format already used by F19:F21 (style "55").
$ws.Range("F19").Copy() | Out-Null
$ws.Range("F22").PasteSpecial(-4122) | Out-Null
$ws.Range("F19").Copy() | Out-Null
$ws.Range("F23").PasteSpecial(-4122) | Out-Null

$ws.Range("F22").Value = 30
$ws.Range("F23").Value = 30
$ws.Range("G22").Formula = "=5.32*2"
$ws.Range("G23").Formula = "=5.32*2"

# N5 is a manually-entered "paid" total that tracks M5 (May); keep it in sync
# with the new SUM(F3:F14) result now that row 7-9 figures changed.
$ws.Range("N5").Value = 278.55

# N6 now folds the new O21 payment into the "paid" total for June.
$ws.Range("N6").Formula = "=SUM(O16,O18,O20:O21)"

# Selection moved to D26 before saving.
$ws.Range("D26").Select() | Out-Null
